# CCS-7, SP-30 Updated proteomics handlers to expect format template from Rolf
#
# The "Value Unit" / "Scale" template on the metadata sheet gains a new
# allowed unit (fmol/ug protein digest), and the data-sheet column headers
# are updated to the new TIME::VALUE_TYPE header format that also encodes
# the biological/technical replicate grouping (::B1_B2::T1_T2).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("openbis-metadata")
$ws2 = $wb.Worksheets.Item("openbis-data")

# --- openbis-metadata: row 6 = "Value Unit" -----------------------------
# Update the description (allowed values) first, then the example value,
# so the shared-string table is rebuilt in the same order Excel used.
$ws1.Range("C6").Value = "One of mM, uM, Percent, RatioT1, RatioCs, or AU, Dimensionless, fmol/ug protein digest"
$ws1.Range("B6").Value = "fmol/ug protein digest"

# --- openbis-data: header row now carries the B1_B2::T1_T2 suffix ------
$ws2.Range("C1").Value = "0::Mean::B1_B2::T1_T2"
$ws2.Range("D1").Value = "0::Std::B1_B2::T1_T2"
# These two headers start with "+", so (as in the real workbook) Excel
# treats them as text via a leading quote, which is why they pick up the
# quotePrefix cell style in the saved file.
$ws2.Range("E1").Value = "'+2100::Mean::B1_B2::T1_T2"
$ws2.Range("F1").Value = "'+2100::Std::B1_B2::T1_T2"

# --- cosmetic: widen columns to fit the new, longer text ---------------
$ws1.Columns("C").AutoFit()
$ws2.Columns("C:F").AutoFit()

# --- cosmetic: restore the cell selections left behind by the edit -----
[void]$ws2.Activate()
[void]$ws2.Range("E6").Select()

[void]$ws1.Activate()
[void]$ws1.Range("C6").Select()
